# Generate Report for handback
# Adds a new handback entry (c44b5bcc-769e-4ed3-9479-ce21f946f2fe) as row 4
# on the "Overview", "zh-cn" and "de-de" worksheets.

$wb = $excel.ActiveWorkbook

$guid    = "c44b5bcc-769e-4ed3-9479-ce21f946f2fe"
$zhFile  = "$guid.5804d24b1de8586c158bb4841f8ed6004f5243f7.zh-cn.xlf"
$deFile  = "$guid.5804d24b1de8586c158bb4841f8ed6004f5243f7.de-de.xlf"
$mdFile  = "$guid.md"
$status  = "Handed back: in sync with en-US"
$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("B4").Value = $status
$wsOverview.Range("C4").Value = $status

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/$mdFile",
    "",
    "",
    $mdFile
)

# ---------------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | Status | Correspond Handoff File |
#                Correspond Handoff Datetime | Target File |
#                Correspond Handback File | Correspond Handback DateTime |
#                Handoff Reason | Dependency From
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B4").Value = $status
$wsZh.Range("D4").Value = "2016-01-28 10:50:26"
$wsZh.Range("D4").NumberFormat = $dateFmt
$wsZh.Range("G4").Value = "2016-01-28 10:51:23"
$wsZh.Range("H4").Value = "Include"

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/$mdFile",
    "",
    "",
    $mdFile
)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhFile",
    "",
    "",
    $zhFile
)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/0000000000000000000000000000000000000000/e2e/$mdFile",
    "",
    "",
    $mdFile
)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0000000000000000000000000000000000000000/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhFile",
    "",
    "",
    $zhFile
)

# ---------------------------------------------------------------------------
# Sheet "de-de": same layout as "zh-cn"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B4").Value = $status
$wsDe.Range("D4").Value = "2016-01-28 10:50:39"
$wsDe.Range("D4").NumberFormat = $dateFmt
$wsDe.Range("G4").Value = "2016-01-28 10:51:45"
$wsDe.Range("H4").Value = "Include"

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/$mdFile",
    "",
    "",
    $mdFile
)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deFile",
    "",
    "",
    $deFile
)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/0000000000000000000000000000000000000000/e2e/$mdFile",
    "",
    "",
    $mdFile
)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0000000000000000000000000000000000000000/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deFile",
    "",
    "",
    $deFile
)
